# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the slide master (the deck's real/active
#                             theme, "Integral" / "Red Violet" colour scheme)
#   ppt/theme/theme2.xml  -> bound to the notes master ("Office Theme" / the
#                             default "Office" colour scheme)
#
# The target edit swaps the two themes' contents: the slide master should end
# up carrying the stock "Office" colour scheme, while the notes master keeps
# the "Red Violet" one. The font scheme and format scheme (fill/line/effect/
# background styles) are identical between the two themes, so the only
# observable difference is the 12-slot colour scheme (dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink).
#
# Re-point the slide master's theme colours at the stock Office palette via
# the supported ThemeColorScheme object model (Item(1..12).RGB, in the
# standard dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order).

$p = $ppt.ActivePresentation
$colors = $p.SlideMaster.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
